# FINFLUX-2815  Stabilaizing automation script
#
# Updates the late-fee % figures (0.23 -> 0.74) and related dependent totals
# across the Summary / Repayment schedule / Transactions sheets, renumbers a
# handful of transaction IDs, and restores the view state (active sheet +
# selected cell per sheet) captured at save time.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A5").Value = 0.74
$wsSummary.Range("E5").Value = 0.74
$wsSummary.Range("F5").Value = 0.74
$wsSummary.Range("B9").Select()

# ---------------------------------------------------------------------
# Repayment schedule sheet
# ---------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Range("J5").Value = 0.74
$wsRepay.Range("K5").Value = 888.46
$wsRepay.Range("Q5").Value = 888.46

# ---------------------------------------------------------------------
# Transactions sheet
# ---------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("A2").Value = 434
$wsTrans.Range("A3").Value = 433
$wsTrans.Range("E3").Value = 23.75
$wsTrans.Range("I3").Value = 0.74
$wsTrans.Range("A4").Value = 425
$wsTrans.Range("A5").Value = 424
$wsTrans.Range("A6").Value = 432
$wsTrans.Range("A7").Value = 423
$wsTrans.Range("D8").Select()

# ---------------------------------------------------------------------
# View state: Repayment schedule ends up the active sheet/tab, with K6
# selected there; activating it last makes it the active tab.
# ---------------------------------------------------------------------
$wsRepay.Activate()
$wsRepay.Range("K6").Select()
